# Edit script: practice task #5 -> #6, and insert "examples" section with a table.

$d = $word.ActiveDocument

# --- 1. "practice task #5: " -> split into 3 runs with "6" instead of "5" ---
$p1 = $d.Paragraphs(1)
$start = $p1.Range.Start
$oldRng = $d.Range($start, $start + 18)
if ($oldRng.Text -eq "practice task #5: ") {
    $oldRng.Text = ""
    $target = $d.Range($start, $start)
    $target.InsertBefore(": ")
    $target.InsertBefore("6")
    $target.InsertBefore("practice task #")
}

# --- 2. Insert "examples" heading + table before the "instructions" heading ---
$instrPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "instructions") {
        $instrPara = $p
        break
    }
}

if ($instrPara -ne $null) {
    $insPos = $instrPara.Range.Start
    $insRng = $d.Range($insPos, $insPos)
    $tablePkg = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="Cmsor2"/></w:pPr><w:r><w:t>examples</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="Rcsostblzat"/><w:tblW w:w="0" w:type="auto"/><w:tblBorders><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tblBorders><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2265"/><w:gridCol w:w="2265"/><w:gridCol w:w="2266"/><w:gridCol w:w="2266"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2265" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">number: 1  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2265" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">guess: 19  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2266" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>score: 0</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2266" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">number: 20  </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2265" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">guess: 18  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2265" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>score: 80</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2266" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">number: 17  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2266" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">guess: 10  </w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2265" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>score: 30</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2265" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">number: 5  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2266" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">guess: 5  </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2266" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>score: 100</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml"><pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="table" w:styleId="Rcsostblzat"><w:name w:val="Table Grid"/><w:basedOn w:val="Normltblzat"/><w:uiPriority w:val="39"/><w:rsid w:val="00CD03AB"/><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:tblPr><w:tblBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tblBorders></w:tblPr></w:style></w:styles></pkg:xmlData></pkg:part></pkg:package>
'@
    $insRng.InsertXML($tablePkg)
}

Write-Output "edit complete"
